# Auto-generated edit script: applies the numeric updates described by the
# Louisoix_Profits.xlsx diff (per-sheet cell value changes, plus a handful of
# cell insertions/deletions where a <c> node was added or removed outright).

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets("ALC")
$ws.Range("H33").Value = 4390
$ws.Range("I33").Value = 5400.579
$ws.Range("J33").Value = 549.8
$ws.Range("K33").Value = 5400.579
$ws.Range("L33").Value = 549.8
$ws.Range("M33").Value = -5171.579
$ws.Range("N33").Value = -1007.8
$ws.Range("H40").Value = 3293.5454
$ws.Range("H64").Value = 23695.715
$ws.Range("J64").Value = 28304.182
$ws.Range("L64").Value = 28304.182
$ws.Range("N64").Value = -28800.182
$ws.Range("H67").Value = 23695.715
$ws.Range("J67").Value = 28304.182
$ws.Range("L67").Value = 28304.182
$ws.Range("N67").Value = -30020.182
$ws.Range("H74").Value = 11151.956
$ws.Range("I74").Value = 7291.923
$ws.Range("J74").Value = 16170
$ws.Range("K74").Value = 7291.923
$ws.Range("L74").Value = 16170
$ws.Range("M74").Value = -6355.923
$ws.Range("N74").Value = -18042
$ws.Range("H77").Value = 11151.956
$ws.Range("I77").Value = 7291.923
$ws.Range("J77").Value = 16170
$ws.Range("K77").Value = 36459.615
$ws.Range("L77").Value = 80850
$ws.Range("M77").Value = -31779.615
$ws.Range("N77").Value = -90210
$ws.Range("H106").Value = 4728.8
$ws.Range("I106").Value = 3111.1667
$ws.Range("K106").Value = 3111.1667
$ws.Range("M106").Value = -2480.1667
$ws.Range("H137").Value = 37152.9
$ws.Range("I137").Value = 42361.92
$ws.Range("K137").Value = 127085.76
$ws.Range("M137").Value = -124535.76
$ws.Range("H138").Value = 2724.7307
$ws.Range("I138").Value = 2954.3572
$ws.Range("J138").Value = 2456.8333
$ws.Range("K138").Value = 8863.071599999999
$ws.Range("L138").Value = 7370.499899999999
$ws.Range("M138").Value = -3723.071599999999
$ws.Range("N138").Value = -17650.4999

# --- Sheet: ARM ---
$ws = $wb.Worksheets("ARM")
$ws.Range("H2").Value = 1089.7778
$ws.Range("J2").Value = 1456.5
$ws.Range("L2").Value = 1456.5
$ws.Range("N2").Value = -1682.5
$ws.Range("H32").Value = 21709
$ws.Range("I32").Value = 25423.357
$ws.Range("J32").Value = 8708.75
$ws.Range("K32").Value = 25423.357
$ws.Range("L32").Value = 8708.75
$ws.Range("M32").Value = -25136.357
$ws.Range("N32").Value = -9282.75
$ws.Range("H116").Value = 1089.7778
$ws.Range("J116").Value = 1456.5
$ws.Range("L116").Value = 1456.5
$ws.Range("N116").Value = -6044.5
$ws.Range("H132").Value = 26993.049
$ws.Range("I132").Value = 32989.453
$ws.Range("J132").Value = 2257.875
$ws.Range("K132").Value = 98968.359
$ws.Range("L132").Value = 6773.625
$ws.Range("M132").Value = -96438.359
$ws.Range("N132").Value = -11833.625

# --- Sheet: BSM ---
$ws = $wb.Worksheets("BSM")
$ws.Range("H3").Value = 1089.7778
$ws.Range("J3").Value = 1456.5
$ws.Range("L3").Value = 1456.5
$ws.Range("N3").Value = -1684.5
$ws.Range("H70").Value = 300000
$ws.Range("J70").Value = 300000
$ws.Range("L70").Value = 300000
$ws.Range("N70").Value = -300586
$ws.Range("H73").Value = 300000
$ws.Range("J73").Value = 300000
$ws.Range("L73").Value = 300000
$ws.Range("N73").Value = -302028

# --- Sheet: CRP ---
$ws = $wb.Worksheets("CRP")
$ws.Range("H31").Value = 2277.6667
$ws.Range("I31").Value = 1653.3846
$ws.Range("J31").Value = 3900.8
$ws.Range("K31").Value = 1653.3846
$ws.Range("L31").Value = 3900.8
$ws.Range("M31").Value = -1358.3846
$ws.Range("N31").Value = -4490.8
$ws.Range("H34").Value = 2277.6667
$ws.Range("I34").Value = 1653.3846
$ws.Range("J34").Value = 3900.8
$ws.Range("K34").Value = 1653.3846
$ws.Range("L34").Value = 3900.8
$ws.Range("M34").Value = -1451.3846
$ws.Range("N34").Value = -4304.8
$ws.Range("H36").Value = 12994.5
$ws.Range("I36").Value = 12994.5
$ws.Range("K36").Value = 12994.5
$ws.Range("M36").Value = -12606.5
$ws.Range("H40").Value = 12994.5
$ws.Range("I40").Value = 12994.5
$ws.Range("K40").Value = 12994.5
$ws.Range("M40").Value = -12834.5
$ws.Range("H42").Value = 9000
$ws.Range("I42").Value = 9000
$ws.Range("K42").Value = 9000
$ws.Range("M42").Value = -8407
$ws.Range("H55").Value = 1999
$ws.Range("I55").Value = 1999
$ws.Range("K55").Value = 1999
$ws.Range("M55").Value = -1684
$ws.Range("H99").Value = 7100
$ws.Range("J99").Value = 2200
$ws.Range("L99").Value = 2200
$ws.Range("N99").Value = -5196
$ws.Range("H122").Value = 3939.8
$ws.Range("I122").Value = 3899
$ws.Range("J122").Value = 3950
$ws.Range("K122").Value = 11697
$ws.Range("L122").Value = 11850
$ws.Range("M122").Value = -9247
$ws.Range("N122").Value = -16750
$ws.Range("H126").Value = 7100
$ws.Range("J126").Value = 2200
$ws.Range("L126").Value = 6600
$ws.Range("N126").Value = -11540

# --- Sheet: CUL ---
$ws = $wb.Worksheets("CUL")
$ws.Range("H12").Value = 113
$ws.Range("J12").Value = 127.55556
$ws.Range("L12").Value = 382.66668
$ws.Range("N12").Value = -728.66668
$ws.Range("H132").Value = 2899
$ws.Range("I132").Value = 2899
$ws.Range("K132").Value = 26091
$ws.Range("M132").Value = -23561

# --- Sheet: GSM ---
$ws = $wb.Worksheets("GSM")
$ws.Range("H52").Value = 32499.334
$ws.Range("I52").Value = 25000
$ws.Range("J52").Value = 39998.668
$ws.Range("K52").Value = 25000
$ws.Range("L52").Value = 39998.668
$ws.Range("M52").Value = -24741
$ws.Range("N52").Value = -40516.668
$ws.Range("H70").Value = 9874.5
$ws.Range("J70").Value = 9874.5
$ws.Range("L70").Value = 9874.5
$ws.Range("N70").Value = -10414.5
$ws.Range("H73").Value = 9874.5
$ws.Range("J73").Value = 9874.5
$ws.Range("L73").Value = 9874.5
$ws.Range("N73").Value = -11746.5
$ws.Range("H97").Value = 3259.6
$ws.Range("I97").Value = 2966.3333
$ws.Range("J97").Value = 3699.5
$ws.Range("K97").Value = 2966.3333
$ws.Range("L97").Value = 3699.5
$ws.Range("M97").Value = -2470.3333
$ws.Range("N97").Value = -4691.5
$ws.Range("H102").Value = 2990.12
$ws.Range("I102").Value = 2879.8823
$ws.Range("K102").Value = 2879.8823
$ws.Range("M102").Value = -1257.8823
$ws.Range("H123").Value = 37997.4
$ws.Range("J123").Value = 37997.4
$ws.Range("L123").Value = 37997.4
$ws.Range("N123").Value = -42897.4
$ws.Range("H126").Value = 6072.4116
$ws.Range("I126").Value = 5787.4546
$ws.Range("J126").Value = 6594.8335
$ws.Range("K126").Value = 17362.3638
$ws.Range("L126").Value = 19784.5005
$ws.Range("M126").Value = -14892.3638
$ws.Range("N126").Value = -24724.5005

# --- Sheet: LTW ---
$ws = $wb.Worksheets("LTW")
$ws.Range("H7").Value = 3227.5715
$ws.Range("I7").Value = 3227.5715
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 3227.5715
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -3115.5715
$ws.Range("N7").ClearContents()
$ws.Range("H16").Value = 8394.799999999999
$ws.Range("J16").Value = 8999
$ws.Range("L16").Value = 8999
$ws.Range("N16").Value = -9339
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()
$ws.Range("H61").Value = 2181.889
$ws.Range("I61").Value = 2140.2856
$ws.Range("K61").Value = 2140.2856
$ws.Range("M61").Value = -1938.2856
$ws.Range("H93").Value = 2857.3333
$ws.Range("I93").Value = 2629.8
$ws.Range("K93").Value = 2629.8
$ws.Range("M93").Value = -1381.8
$ws.Range("H113").Value = 2181.889
$ws.Range("I113").Value = 2140.2856
$ws.Range("K113").Value = 2140.2856
$ws.Range("M113").Value = 29.71439999999984
$ws.Range("H122").Value = 3297.75
$ws.Range("I122").Value = 2982.1667
$ws.Range("K122").Value = 8946.500100000001
$ws.Range("M122").Value = -6496.500100000001
$ws.Range("H126").Value = 3227.5715
$ws.Range("I126").Value = 3227.5715
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 9682.7145
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -7212.7145
$ws.Range("N126").ClearContents()

# --- Sheet: WVR ---
$ws = $wb.Worksheets("WVR")
$ws.Range("H55").Value = 9666.333000000001
$ws.Range("I55").Value = 4499.5
$ws.Range("J55").Value = 20000
$ws.Range("K55").Value = 4499.5
$ws.Range("L55").Value = 20000
$ws.Range("M55").Value = -4222.5
$ws.Range("N55").Value = -20554
$ws.Range("H122").Value = 2176.9375
$ws.Range("I122").Value = 1959.0344
$ws.Range("K122").Value = 5877.1032
$ws.Range("M122").Value = -3427.1032
$ws.Range("H126").Value = 39172.332
$ws.Range("I126").Value = 46996.453
$ws.Range("J126").Value = 4746.2
$ws.Range("K126").Value = 140989.359
$ws.Range("L126").Value = 14238.6
$ws.Range("M126").Value = -138519.359
$ws.Range("N126").Value = -19178.6
